# Automatic map update (mapa_interactivo.html)
# A new field record (Caso -500, "Castañares 5656") replaces the previous
# last record (Caso -515, "Rivadavia 7470") in the dataset. Because the
# data is ordered chronologically (column B) rather than by Caso id, the
# new record is inserted in the middle of the "General" master list (and
# of the matching "AYKO" provider-filtered list) and every following row
# shifts down by one; the row that drops off the bottom (the old -515
# "Rivadavia 7470" record, now duplicated after the shift) is deleted.

$wb = $excel.ActiveWorkbook

function Insert-NewCaseRow {
    param(
        [string]$SheetName,
        [int]$InsertAtRow,
        [int]$OldLastRow
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Shift InsertAtRow..OldLastRow down by one, opening a blank row.
    $ws.Rows.Item($InsertAtRow).Insert()

    # Populate the newly opened row with the new case data.
    $ws.Cells.Item($InsertAtRow, 1).Value = "'-500"
    $ws.Cells.Item($InsertAtRow, 2).Value = "'7/3/2025"
    $ws.Cells.Item($InsertAtRow, 3).Value = "Castañares 5656"
    $ws.Cells.Item($InsertAtRow, 4).Value = "'8"
    $ws.Cells.Item($InsertAtRow, 5).Value = "'807965768"
    $ws.Cells.Item($InsertAtRow, 6).Value = "AYKO"
    $ws.Cells.Item($InsertAtRow, 7).Value = "Pendiente"
    $ws.Cells.Item($InsertAtRow, 8).Value = "Columna chocada con rienda a pique"
    $ws.Cells.Item($InsertAtRow, 9).Value = 1
    $ws.Cells.Item($InsertAtRow, 10).Value = "Cambio"
    $ws.Cells.Item($InsertAtRow, 11).Value = "Sin equipos"
    $ws.Cells.Item($InsertAtRow, 12).Value = "Terminal"
    $ws.Cells.Item($InsertAtRow, 13).Value = -58.479921
    $ws.Cells.Item($InsertAtRow, 14).Value = -34.673021
    $ws.Cells.Item($InsertAtRow, 15).Value = "Boedo"
    $ws.Cells.Item($InsertAtRow, 16).Value = "Capital Sur"

    # Every later row (including the old last row) has now shifted down by
    # one, so the row that used to be OldLastRow now duplicates itself at
    # OldLastRow + 1. Remove that trailing duplicate so the sheet keeps its
    # original row count.
    $ws.Rows.Item($OldLastRow + 1).Delete()
}

# "General" master sheet: new row becomes row 384, old data 384-412 shifts
# to 385-413, and the old row 412 (now duplicated at 413) is removed.
Insert-NewCaseRow "General" 384 412

# "AYKO" provider sheet: same record, same shift, local row numbers.
Insert-NewCaseRow "AYKO" 89 91
